# Word COM-interop edit script
# Applies:
#  1) Split the "<emotion>_<image name>.<image type>" sentence into three
#     runs, wrapping the "&gt;.&lt;" fragment with proofErr gramStart/gramEnd.
#  2) Split the "... predicted by dlib." sentence into three runs, wrapping
#     "dlib" with proofErr spellStart/spellEnd.
#  3) Move the "_GoBack" bookmark from the "Instead o|f predicting" split
#     up to sit right after the "Conclusion" heading run.

$d = $word.ActiveDocument

# Pull the full WordprocessingML for the main document story so we can
# perform precise, surgical text/run surgery (Word's object model doesn't
# give a direct "split this run and stick a proofErr between the pieces"
# verb, so we edit the underlying markup for the Word.Document directly,
# the same markup Find/Replace and Range edits above would otherwise
# produce piecemeal).
$xml = $d.WordOpenXML

# --- 1) "<emotion>_<image name>.<image type>" -------------------------
$old1 = '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>The category for each image is specified before the image name. The format is &lt;emotion&gt;_&lt;image name&gt;.&lt;image type&gt;</w:t></w:r>'
$new1 = '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>The category for each image is specified before the image name. The format is &lt;emotion&gt;_&lt;image name</w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>&gt;.&lt;</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>image type&gt;</w:t></w:r>'

if ($xml.IndexOf($old1) -lt 0) {
    throw "Could not locate the '<emotion>_<image name>...' run to split"
}
$xml = $xml.Replace($old1, $new1)

# --- 2) "... predicted by dlib." ---------------------------------------
$old2 = '<w:r w:rsidRPr="00E6546B"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>All the models are fed in the x and y position for each of the 68 landmarks predicted by dlib.</w:t></w:r>'
$new2 = '<w:r w:rsidRPr="00E6546B"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">All the models are fed in the x and y position for each of the 68 landmarks predicted by </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>dlib</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.</w:t></w:r>'

if ($xml.IndexOf($old2) -lt 0) {
    throw "Could not locate the 'predicted by dlib.' run to split"
}
$xml = $xml.Replace($old2, $new2)

# --- 3) Move the _GoBack bookmark up to the "Conclusion" heading -------
$old3 = '<w:bookmarkStart w:id="7" w:name="_Toc37081444"/><w:r><w:lastRenderedPageBreak/><w:t>Conclusion</w:t></w:r><w:bookmarkEnd w:id="7"/>'
$new3 = '<w:bookmarkStart w:id="7" w:name="_Toc37081444"/><w:r><w:lastRenderedPageBreak/><w:t>Conclusion</w:t></w:r><w:bookmarkStart w:id="8" w:name="_GoBack"/><w:bookmarkEnd w:id="7"/><w:bookmarkEnd w:id="8"/>'

if ($xml.IndexOf($old3) -lt 0) {
    throw "Could not locate the Conclusion heading bookmark"
}
$xml = $xml.Replace($old3, $new3)

$old4 = '<w:r w:rsidR="00E6546B"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>o</w:t></w:r><w:bookmarkStart w:id="8" w:name="_GoBack"/><w:bookmarkEnd w:id="8"/><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>f predicting'
$new4 = '<w:r w:rsidR="00E6546B"><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>o</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>f predicting'

if ($xml.IndexOf($old4) -lt 0) {
    throw "Could not locate the old _GoBack bookmark position"
}
$xml = $xml.Replace($old4, $new4)

# Write the rebuilt markup back to the document.
$d.WordOpenXML = $xml

Write-Host "Edits applied."
